$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column I (9). This shifts old column I -> J,
#    leaves B..H untouched, creates a new blank column I whose per-row cell
#    styles are copied from the former column H, and widens the B2:I2 merge
#    to B2:J2 automatically.
$ws.Columns.Item(9).Insert()

# 2. The "Generated Date :" / "Generated By :" label+value pair that used to
#    live in H5:I5 and H6:I6 actually belongs one column further right now
#    (I5:J5 / I6:J6), leaving the newly inserted column H blank on those two
#    rows. Move the values/styles over by hand (cut/paste across an
#    overlapping range is unreliable, so read, clear, then write).
$h5Value = $ws.Range("H5").Value2
$h6Value = $ws.Range("H6").Value2

$ws.Range("I5").Value = $h5Value
$ws.Range("I5").Font.Name = "Times New Roman"
$ws.Range("I5").Font.Size = 12
$ws.Range("I5").Font.Bold = $true
$ws.Range("I5").HorizontalAlignment = -4131
$ws.Range("I5").VerticalAlignment = -4108

$ws.Range("I6").Value = $h6Value
$ws.Range("I6").Font.Name = "Times New Roman"
$ws.Range("I6").Font.Size = 12
$ws.Range("I6").Font.Bold = $true
$ws.Range("I6").HorizontalAlignment = -4131
$ws.Range("I6").VerticalAlignment = -4108

$ws.Range("H5").Clear()
$ws.Range("H5").NumberFormat = "General"

$ws.Range("H6").Clear()
$ws.Range("H6").NumberFormat = "General"

# J6 becomes the (currently empty) value cell for "Generated By :" and needs
# its own left/vertical-centered, non-bold text style (distinct from the
# date-valued J5 cell).
$ws.Range("J6").Font.Name = "Times New Roman"
$ws.Range("J6").Font.Size = 12
$ws.Range("J6").Font.Bold = $false
$ws.Range("J6").HorizontalAlignment = -4131
$ws.Range("J6").VerticalAlignment = -4108
$ws.Range("J6").NumberFormat = "General"

# 3. New "Created By" header label, added to the table header row (row 8),
#    in the newly inserted column I; the old header ("Date Created") that
#    used to sit there is now in J8 (already shifted there by the column
#    insert in step 1).
$ws.Range("I8").Value = "Created By"

# 4. Fix up the column widths: the newly inserted column H gets a new,
#    wider width, while column I inherits the old column H's width/format.
$ws.Columns.Item(8).ColumnWidth = 19.5
$ws.Columns.Item(9).ColumnWidth = 16.666666666666668
$ws.Columns.Item(9).HorizontalAlignment = -4108

# 5. Touch A1 (as a true visual no-op: border style is already "none") so the
#    sheet's used range / dimension starts at column A again, matching the
#    template's originally authored dimension.
$ws.Range("A1").Borders.LineStyle = -4142

# 6. Restore the original active-cell selection recorded in the file.
$ws.Range("J15").Select()
